# Add two new LeetCode tracker rows (#74 "Search a 2D Matrix" and
# #162 "Find Peak Element") to the bottom of the log table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 24: LeetCode 74, Search a 2D Matrix ---------------------------
$ws.Cells.Item(24, 1).Value = 74
$ws.Cells.Item(24, 2).Value = "Search a 2D Matrix"
$ws.Cells.Item(24, 3).Value = "#array  #binary-search #matrix #核心 "
$ws.Cells.Item(24, 4).Value = "medium"
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 10

# Copy the date formatting from the row above (style 3 = short-date,
# centered) before writing the serial date values so H24/I24 pick up
# the same number format as every other "Time"/"Last Update" cell.
$ws.Cells.Item(23, 8).Copy()
$ws.Cells.Item(24, 8).PasteSpecial(-4122)
$ws.Cells.Item(24, 8).Value = 45838

$ws.Cells.Item(23, 9).Copy()
$ws.Cells.Item(24, 9).PasteSpecial(-4122)
$ws.Cells.Item(24, 9).Value = 45838

$ws.Rows.Item(24).RowHeight = 51

# --- Row 25: LeetCode 162, Find Peak Element ----------------------------
$ws.Cells.Item(25, 1).Value = 162
$ws.Cells.Item(25, 2).Value = "Find Peak Element"
$ws.Cells.Item(25, 3).Value = "#array #binary-search #核心 "
$ws.Cells.Item(25, 4).Value = "medium"
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 18

$ws.Cells.Item(23, 8).Copy()
$ws.Cells.Item(25, 8).PasteSpecial(-4122)
$ws.Cells.Item(25, 8).Value = 45838

$ws.Cells.Item(23, 9).Copy()
$ws.Cells.Item(25, 9).PasteSpecial(-4122)
$ws.Cells.Item(25, 9).Value = 45838

$ws.Rows.Item(25).RowHeight = 34

# Move the view/selection down to the newly added last row, matching the
# author's final cursor position in the sheet.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("I25").Select()
